$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per corrected naive component forecaster calculation
$ws.Range("J25").Value = -4.341321906298077
$ws.Range("K25").Value = 0.5976254839025203
$ws.Range("I26").Value = -4.365333711727925
$ws.Range("J26").Value = 0.5672256391354182
$ws.Range("H27").Value = -4.285655806507345
$ws.Range("I27").Value = 0.6532326549931144
$ws.Range("G28").Value = -4.265333711727925
$ws.Range("H28").Value = 0.6672256391354182
$ws.Range("F29").Value = -4.284477149457854
$ws.Range("G29").Value = 0.6560769697073342
$ws.Range("H29").Value = 1.658114585577847
$ws.Range("I29").Value = -1.966123838979201
$ws.Range("J29").Value = -1.647175090369943
$ws.Range("K29").Value = 0.9415002252234601
$ws.Range("E30").Value = -4.305330271545129
$ws.Range("F30").Value = 0.6446265913711419
$ws.Range("G30").Value = 1.761846783446432
$ws.Range("H30").Value = -1.912723000994306
$ws.Range("I30").Value = -1.614292410696496
$ws.Range("J30").Value = 0.9300355966686533
$ws.Range("D31").Value = -2.512996718662795
$ws.Range("E31").Value = 1.748381033365794
$ws.Range("F31").Value = 2.917126296654925
$ws.Range("G31").Value = -1.273578130868497
$ws.Range("H31").Value = -1.611081864576003
$ws.Range("I31").Value = 0.9331061597207662
$ws.Range("C32").Value = -0.165333711727925
$ws.Range("D32").Value = 1.667225639135418
$ws.Range("E32").Value = 2.46956845314902
$ws.Range("F32").Value = -1.250117300589673
$ws.Range("G32").Value = -1.531122027629067
$ws.Range("H32").Value = 0.9497577098212645
$ws.Range("B33").Value = 0.06406792636339276
$ws.Range("C33").Value = 2.154138501756804
$ws.Range("D33").Value = 2.059333715820543
$ws.Range("E33").Value = -1.619557144118616
$ws.Range("F33").Value = -1.609340397885774
$ws.Range("G33").Value = 1.000755198365638
$ws.Range("H33").Value = -1.324862745064436
$ws.Range("I33").Value = -0.262338436028287
$ws.Range("J33").Value = -0.07586005108070992
$ws.Range("K33").Value = 1.015792704729378
$ws.Range("B34").Value = -0.6232870972609703
$ws.Range("C34").Value = 2.159983717119644
$ws.Range("D34").Value = 0.145870948706488
$ws.Range("E34").Value = -1.159709944833494
$ws.Range("F34").Value = 1.254082278054852
$ws.Range("G34").Value = -1.291125139570511
$ws.Range("H34").Value = -0.2417684416796248
$ws.Range("I34").Value = -0.05767949661618843
$ws.Range("J34").Value = 1.025133804506055
$ws.Range("B35").Value = -0.77312532813319
$ws.Range("C35").Value = -0.01313334549770717
$ws.Range("D35").Value = 0.711742572249292
$ws.Range("E35").Value = 1.765408493094085
$ws.Range("F35").Value = -1.27419524477169
$ws.Range("G35").Value = -0.2210300370960283
$ws.Range("H35").Value = -0.01870964234420563
$ws.Range("I35").Value = 1.02969361131457
$ws.Range("B36").Value = -0.2570926644240221
$ws.Range("C36").Value = 2.360555864296444
$ws.Range("D36").Value = 2.132036102469129
$ws.Range("E36").Value = -1.302806821653718
$ws.Range("F36").Value = -0.2121505183469736
$ws.Range("G36").Value = -0.01567838362100815
$ws.Range("H36").Value = 1.036859485198733
$ws.Range("B37").Value = -0.9325864806726689
$ws.Range("C37").Value = 2.176039623185105
$ws.Range("D37").Value = -0.283047497385382
$ws.Range("E37").Value = 0.3346930603001257
$ws.Range("F37").Value = 0.2620123033535435
$ws.Range("G37").Value = 1.083095960590256
$ws.Range("H37").Value = 1.26469213391276
$ws.Range("I37").Value = -0.3325428966199979
$ws.Range("J37").Value = 0.5727158612572509
$ws.Range("K37").Value = 0.4929530572011345
$ws.Range("B38").Value = 0.07804016256402241
$ws.Range("C38").Value = -1.997727010608287
$ws.Range("D38").Value = 1.50224383099669
$ws.Range("E38").Value = 1.088925148385484
$ws.Range("F38").Value = 1.208640679979695
$ws.Range("G38").Value = 1.419101529386267
$ws.Range("H38").Value = -0.2565607374317693
$ws.Range("I38").Value = 0.627946269261713
$ws.Range("J38").Value = 0.5482123265480681
$ws.Range("B39").Value = -0.8350711388588363
$ws.Range("C39").Value = 0.1485789965021382
$ws.Range("D39").Value = 0.3742447754406868
$ws.Range("E39").Value = 1.537655471550067
$ws.Range("F39").Value = 1.821898845291485
$ws.Range("G39").Value = 0.1769899521501799
$ws.Range("H39").Value = 0.8608096660974809
$ws.Range("I39").Value = 0.6381556253882366
$ws.Range("B40").Value = -0.0806700814468968
$ws.Range("C40").Value = 0.4084836957640618
$ws.Range("D40").Value = 1.497444495070027
$ws.Range("E40").Value = 1.724547918235458
$ws.Range("F40").Value = 0.178741787876163
$ws.Range("G40").Value = 1.003530935625179
$ws.Range("H40").Value = 0.7198339373236422
$ws.Range("B41").Value = -0.43807954759518
$ws.Range("C41").Value = 0.1477304037695859
$ws.Range("D41").Value = 0.4534906893407109
$ws.Range("E41").Value = -0.5222947029782006
$ws.Range("F41").Value = 0.7510341514886107
$ws.Range("G41").Value = 0.8765966212947234
$ws.Range("H41").Value = 0.6504644380763409
$ws.Range("I41").Value = 0.6588952181776051
$ws.Range("J41").Value = -0.4009133581774047
$ws.Range("K41").Value = -0.3987502414539961
$ws.Range("B42").Value = -0.4230865868247518
$ws.Range("C42").Value = 0.4652481907266262
$ws.Range("D42").Value = -0.5633934991668781
$ws.Range("E42").Value = 0.9107937393168015
$ws.Range("F42").Value = 0.8263562670369282
$ws.Range("G42").Value = 0.5509036232411972
$ws.Range("H42").Value = 0.4448453594173571
$ws.Range("I42").Value = -0.640932475279683
$ws.Range("J42").Value = -0.6738259216274827
$ws.Range("B43").Value = 0.4354409026540649
$ws.Range("C43").Value = -0.4697653985074481
$ws.Range("D43").Value = 0.8383240611323401
$ws.Range("E43").Value = 0.8745247415661055
$ws.Range("F43").Value = 0.5977696827513626
$ws.Range("G43").Value = 0.5182375512410431
$ws.Range("H43").Value = -0.6131448997133897
$ws.Range("I43").Value = -0.7044789358240422
$ws.Range("B44").Value = -0.3395690612336324
$ws.Range("C44").Value = 0.5036576038034468
$ws.Range("D44").Value = 0.7113494639342406
$ws.Range("E44").Value = 0.4670600820553261
$ws.Range("F44").Value = 0.4344288604015389
$ws.Range("G44").Value = -0.6674920639201498
$ws.Range("H44").Value = -0.7053026745305289
$ws.Range("B45").Value = 0.3275174550509519
$ws.Range("C45").Value = 0.6384202572023839
$ws.Range("D45").Value = 0.426465400382952
$ws.Range("E45").Value = 0.471006240264856
$ws.Range("F45").Value = -0.5501478897554928
$ws.Range("G45").Value = -0.6532877491573856
$ws.Range("H45").Value = 0.9228462089976119
$ws.Range("I45").Value = -0.1008087126355974
$ws.Range("B46").Value = 0.2089562936243113
$ws.Range("C46").Value = 0.1852825741506052
$ws.Range("D46").Value = 0.3176569333398902
$ws.Range("E46").Value = -0.5556119317386812
$ws.Range("F46").Value = -0.5290948250225114
$ws.Range("G46").Value = 0.9965123179614324
$ws.Range("H46").Value = -0.05555969787798182
$ws.Range("B47").Value = 0.02912383308249389
$ws.Range("C47").Value = 0.2041938464851824
$ws.Range("D47").Value = -0.5155302980697907
$ws.Range("E47").Value = -0.4334147766432039
$ws.Range("F47").Value = 1.047431584022505
$ws.Range("G47").Value = -0.001067671251348712
$ws.Range("B48").Value = -0.1406678742931149
$ws.Range("C48").Value = -0.7175887441063646
$ws.Range("D48").Value = -0.5795609912632642
$ws.Range("E48").Value = 1.01571912334704
$ws.Range("F48").Value = 0.1211914621108861
$ws.Range("B49").Value = -1.067943258854512
$ws.Range("C49").Value = -0.7665539717803824
$ws.Range("D49").Value = 0.8995656217062589
$ws.Range("E49").Value = 0.05336612861586332
$ws.Range("B50").Value = -0.8299273031874748
$ws.Range("C50").Value = 0.7599009234537135
$ws.Range("D50").Value = -0.2351729235711133
$ws.Range("B51").Value = 1.026563613664763
$ws.Range("C51").Value = -0.221020254317267
$ws.Range("B52").Value = -0.2176617297482864

# Clear cells that are no longer part of the (now one-column-shorter) data block
$ws.Range("J45").ClearContents()
$ws.Range("I46").ClearContents()
$ws.Range("H47").ClearContents()
$ws.Range("G48").ClearContents()
$ws.Range("F49").ClearContents()
$ws.Range("E50").ClearContents()
$ws.Range("D51").ClearContents()
$ws.Range("C52").ClearContents()
$ws.Range("B53").ClearContents()
